$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A-E and G hold plain text in the source data (dates spelled out,
# day names, time-of-day strings, ids, free text, file names), and column F
# ("confidence") holds numbers formatted as plain text (e.g. "60.07") in the
# original file. Force text number format on the whole block we're about to
# (re)write so Excel doesn't reinterpret any of these as real numbers/dates,
# matching the original inlineStr layout - then restore the default "Normal"
# style afterwards so no stray cell style gets introduced.
$dataRange = $ws.Range("A2:G7")
$dataRange.NumberFormat = "@"

# Row 2
$ws.Cells.Item(2,1).Value = "05 July 2025"
$ws.Cells.Item(2,2).Value = "Saturday"
$ws.Cells.Item(2,3).Value = "12:17:02"
$ws.Cells.Item(2,4).Value = "M16"
$ws.Cells.Item(2,5).Value = "mobil"
$ws.Cells.Item(2,6).Value = "60.07"
$ws.Cells.Item(2,7).Value = "pelanggaran_20250705_121702.jpg"

# Row 3
$ws.Cells.Item(3,1).Value = "05 July 2025"
$ws.Cells.Item(3,2).Value = "Saturday"
$ws.Cells.Item(3,3).Value = "12:17:02"
$ws.Cells.Item(3,4).Value = "M17"
$ws.Cells.Item(3,5).Value = "mobil"
$ws.Cells.Item(3,6).Value = "60.07"
$ws.Cells.Item(3,7).Value = "pelanggaran_20250705_121702.jpg"

# Row 4
$ws.Cells.Item(4,1).Value = "05 July 2025"
$ws.Cells.Item(4,2).Value = "Saturday"
$ws.Cells.Item(4,3).Value = "12:17:02"
$ws.Cells.Item(4,4).Value = "M19"
$ws.Cells.Item(4,5).Value = "mobil"
$ws.Cells.Item(4,6).Value = "60.07"
$ws.Cells.Item(4,7).Value = "pelanggaran_20250705_121702.jpg"

# Row 5 (new)
$ws.Cells.Item(5,1).Value = "05 July 2025"
$ws.Cells.Item(5,2).Value = "Saturday"
$ws.Cells.Item(5,3).Value = "12:17:02"
$ws.Cells.Item(5,4).Value = "M23"
$ws.Cells.Item(5,5).Value = "mobil"
$ws.Cells.Item(5,6).Value = "60.07"
$ws.Cells.Item(5,7).Value = "pelanggaran_20250705_121702.jpg"

# Row 6 (new)
$ws.Cells.Item(6,1).Value = "05 July 2025"
$ws.Cells.Item(6,2).Value = "Saturday"
$ws.Cells.Item(6,3).Value = "12:24:54"
$ws.Cells.Item(6,4).Value = "M24"
$ws.Cells.Item(6,5).Value = "mobil"
$ws.Cells.Item(6,6).Value = "79.27"
$ws.Cells.Item(6,7).Value = "pelanggaran_20250705_122454.jpg"

# Row 7 (new)
$ws.Cells.Item(7,1).Value = "05 July 2025"
$ws.Cells.Item(7,2).Value = "Saturday"
$ws.Cells.Item(7,3).Value = "13:43:56"
$ws.Cells.Item(7,4).Value = "M112"
$ws.Cells.Item(7,5).Value = "mobil"
$ws.Cells.Item(7,6).Value = "42.60"
$ws.Cells.Item(7,7).Value = "pelanggaran_20250705_134356.jpg"

# Restore the default "Normal" cell style so no extra style index is
# introduced by the temporary text-number-format override above.
$dataRange.Style = "Normal"

# Update the sheet view's selection to match the new active cell (diff:
# H17 -> H11).
$ws.Range("H11").Select()
